{"js": "// Replace the header date and each two-digit multiplication answer cell\n// with the values from the target revision. Every \"old\" string below is\n// unique within the document, so a direct search+replace is unambiguous.\nconst replacements = [\n  [\"2025-12-02 Tuesday\", \"2025-12-03 Wednesday\"],\n  [\"17\u00d753=901\", \"59\u00d745=2655\"],\n  [\"77\u00d774=5698\", \"51\u00d738=1938\"],\n  [\"70\u00d770=4900\", \"85\u00d736=3060\"],\n  [\"59\u00d764=3776\", \"91\u00d749=4459\"],\n  [\"84\u00d715=1260\", \"52\u00d797=5044\"],\n  [\"49\u00d786=4214\", \"20\u00d722=440\"],\n  [\"95\u00d725=2375\", \"72\u00d713=936\"],\n  [\"18\u00d750=900\", \"64\u00d746=2944\"],\n  [\"64\u00d751=3264\", \"85\u00d724=2040\"],\n  [\"67\u00d715=1005\", \"61\u00d798=5978\"],\n  [\"64\u00d741=2624\", \"20\u00d723=460\"],\n  [\"15\u00d740=600\", \"28\u00d793=2604\"],\n  [\"95\u00d774=7030\", \"16\u00d786=1376\"],\n  [\"31\u00d737=1147\", \"94\u00d785=7990\"],\n  [\"58\u00d718=1044\", \"27\u00d769=1863\"],\n  [\"36\u00d718=648\", \"81\u00d775=6075\"],\n  [\"42\u00d795=3990\", \"71\u00d755=3905\"],\n  [\"17\u00d725=425\", \"81\u00d773=5913\"],\n  [\"49\u00d758=2842\", \"27\u00d719=513\"],\n  [\"60\u00d735=2100\", \"33\u00d793=3069\"],\n  [\"64\u00d728=1792\", \"13\u00d754=702\"],\n  [\"74\u00d714=1036\", \"17\u00d735=595\"],\n  [\"21\u00d788=1848\", \"34\u00d728=952\"],\n  [\"80\u00d727=2160\", \"49\u00d739=1911\"],\n  [\"80\u00d790=7200\", \"15\u00d745=675\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the header date and each two-digit multiplication answer cell to\n# match the target revision. Every \"Old\" value below occurs exactly once in\n# the document, so Find/Replace across the whole body is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @{ Old = \"2025-12-02 Tuesday\"; New = \"2025-12-03 Wednesday\" },\n  @{ Old = \"17\u00d753=901\"; New = \"59\u00d745=2655\" },\n  @{ Old = \"77\u00d774=5698\"; New = \"51\u00d738=1938\" },\n  @{ Old = \"70\u00d770=4900\"; New = \"85\u00d736=3060\" },\n  @{ Old = \"59\u00d764=3776\"; New = \"91\u00d749=4459\" },\n  @{ Old = \"84\u00d715=1260\"; New = \"52\u00d797=5044\" },\n  @{ Old = \"49\u00d786=4214\"; New = \"20\u00d722=440\" },\n  @{ Old = \"95\u00d725=2375\"; New = \"72\u00d713=936\" },\n  @{ Old = \"18\u00d750=900\"; New = \"64\u00d746=2944\" },\n  @{ Old = \"64\u00d751=3264\"; New = \"85\u00d724=2040\" },\n  @{ Old = \"67\u00d715=1005\"; New = \"61\u00d798=5978\" },\n  @{ Old = \"64\u00d741=2624\"; New = \"20\u00d723=460\" },\n  @{ Old = \"15\u00d740=600\"; New = \"28\u00d793=2604\" },\n  @{ Old = \"95\u00d774=7030\"; New = \"16\u00d786=1376\" },\n  @{ Old = \"31\u00d737=1147\"; New = \"94\u00d785=7990\" },\n  @{ Old = \"58\u00d718=1044\"; New = \"27\u00d769=1863\" },\n  @{ Old = \"36\u00d718=648\"; New = \"81\u00d775=6075\" },\n  @{ Old = \"42\u00d795=3990\"; New = \"71\u00d755=3905\" },\n  @{ Old = \"17\u00d725=425\"; New = \"81\u00d773=5913\" },\n  @{ Old = \"49\u00d758=2842\"; New = \"27\u00d719=513\" },\n  @{ Old = \"60\u00d735=2100\"; New = \"33\u00d793=3069\" },\n  @{ Old = \"64\u00d728=1792\"; New = \"13\u00d754=702\" },\n  @{ Old = \"74\u00d714=1036\"; New = \"17\u00d735=595\" },\n  @{ Old = \"21\u00d788=1848\"; New = \"34\u00d728=952\" },\n  @{ Old = \"80\u00d727=2160\"; New = \"49\u00d739=1911\" },\n  @{ Old = \"80\u00d790=7200\"; New = \"15\u00d745=675\" }\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $found = $range.Find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n    if (-not $found) {\n        throw \"Text not found: $($pair.Old)\"\n    }\n}\n"}
